# Apply the "june - oct 2025" update to the shuttlecock purchase workbook.
# Adds a new purchase record (row 21) to the "shuttlecock_buy" sheet,
# extends the shared formulas for columns G, H and I down to row 21,
# and updates the active selection to the newly added cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shuttlecock_buy")
$ws.Activate()

# New row of data (index 20 in column A, continuing the existing sequence).
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 45962
$ws.Range("B21").NumberFormat = $ws.Range("B20").NumberFormat

$ws.Range("C21").Value = "Average"

$ws.Range("D21").Value = 720
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0

# Extend the shared formulas from row 20 down into the new row 21.
$ws.Range("G21").Formula = "=D21*E21+F21"
$ws.Range("H21").Formula = "=E21*12"
$ws.Range("I21").Formula = "=ROUNDUP(G21/H21,0)"

# Update the selection to match the final state of the sheet.
$ws.Range("G21").Select()

$wb.Save()
